$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176:274 down to 177:275
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new data record
$ws.Cells.Item(176, 1).Value = 10
$ws.Cells.Item(176, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(176, 3).Value = 'La Araucanía'
$ws.Cells.Item(176, 4).Value = 44606
$ws.Cells.Item(176, 5).Value = 9
$ws.Cells.Item(176, 6).Value = 100112009
$ws.Cells.Item(176, 7).Value = 'Acelga'
$ws.Cells.Item(176, 8).Value = 'Sin especificar'
$ws.Cells.Item(176, 9).Value = 'Primera'
$ws.Cells.Item(176, 10).Value = 80
$ws.Cells.Item(176, 11).Value = 8000
$ws.Cells.Item(176, 12).Value = 8000
$ws.Cells.Item(176, 13).Value = 8000
$ws.Cells.Item(176, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(176, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(176, 16).Value = 667
$ws.Cells.Item(176, 17).Value = 12
$ws.Cells.Item(176, 18).Value = 'Hortaliza'
